# Update the division-problem table cells to the new values.
# Cells are addressed by (row, column) in document order since some
# source values (e.g. "91÷7=") repeat with different replacements.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{Row=1;  Col=1; Old="91÷7="; New="58÷4="},
    @{Row=1;  Col=2; Old="65÷7="; New="51÷7="},
    @{Row=1;  Col=3; Old="60÷5="; New="23÷7="},
    @{Row=1;  Col=4; Old="91÷4="; New="94÷6="},
    @{Row=1;  Col=5; Old="57÷3="; New="54÷8="},

    @{Row=5;  Col=1; Old="32÷7="; New="23÷3="},
    @{Row=5;  Col=2; Old="64÷8="; New="98÷6="},
    @{Row=5;  Col=3; Old="17÷7="; New="18÷7="},
    @{Row=5;  Col=4; Old="66÷2="; New="37÷3="},
    @{Row=5;  Col=5; Old="64÷9="; New="65÷6="},

    @{Row=9;  Col=1; Old="99÷8="; New="70÷8="},
    @{Row=9;  Col=2; Old="55÷8="; New="88÷6="},
    @{Row=9;  Col=3; Old="77÷4="; New="75÷6="},
    @{Row=9;  Col=4; Old="29÷2="; New="30÷6="},
    @{Row=9;  Col=5; Old="95÷2="; New="54÷7="},

    @{Row=13; Col=1; Old="78÷7="; New="35÷5="},
    @{Row=13; Col=2; Old="61÷5="; New="80÷9="},
    @{Row=13; Col=3; Old="34÷2="; New="53÷8="},
    @{Row=13; Col=4; Old="43÷5="; New="95÷9="},
    @{Row=13; Col=5; Old="49÷8="; New="24÷7="},

    @{Row=17; Col=1; Old="77÷5="; New="72÷4="},
    @{Row=17; Col=2; Old="80÷6="; New="80÷7="},
    @{Row=17; Col=3; Old="17÷3="; New="97÷3="},
    @{Row=17; Col=4; Old="82÷2="; New="89÷3="},
    @{Row=17; Col=5; Old="91÷7="; New="77÷2="}
)

foreach ($c in $changes) {
    $cell = $t.Cell($c.Row, $c.Col)
    $cell.Range.Text = $c.New
}
